$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")
$ws.Cells.Item(2, 1).Value = "30 Oct 2025, 12:46 PM"

$ws = $wb.Worksheets.Item("Top Gainers")
$ws.Cells.Item(2, 3).Value = 11.4967
$ws.Cells.Item(2, 4).Value = 20.3284
$ws.Cells.Item(2, 5).Value = 27.5068
$ws.Cells.Item(3, 3).Value = 10.8263
$ws.Cells.Item(3, 4).Value = 11.0334
$ws.Cells.Item(3, 5).Value = 25.1404
$ws.Cells.Item(4, 3).Value = 10.2538
$ws.Cells.Item(4, 4).Value = 15.5319
$ws.Cells.Item(4, 5).Value = 22.4352
$ws.Cells.Item(5, 3).Value = 9.3965
$ws.Cells.Item(5, 4).Value = 17.1086
$ws.Cells.Item(5, 5).Value = 30.4501
$ws.Cells.Item(6, 3).Value = 9.041
$ws.Cells.Item(6, 4).Value = 6.1254
$ws.Cells.Item(6, 5).Value = -8.5463
$ws.Cells.Item(7, 3).Value = 8.3438
$ws.Cells.Item(7, 4).Value = 10.6342
$ws.Cells.Item(7, 5).Value = 26.9236
$ws.Cells.Item(8, 3).Value = 7.5584
$ws.Cells.Item(8, 4).Value = 12.7598
$ws.Cells.Item(8, 5).Value = 14.9248
$ws.Cells.Item(9, 3).Value = 7.1871
$ws.Cells.Item(9, 4).Value = 3.6886
$ws.Cells.Item(9, 5).Value = 16.3506
$ws.Cells.Item(10, 2).Value = "SHREEJISPG"
$ws.Cells.Item(10, 3).Value = 5.9962
$ws.Cells.Item(10, 4).Value = 10.2659
$ws.Cells.Item(10, 5).Value = 11.3595
$ws.Cells.Item(11, 2).Value = "MCLOUD"
$ws.Cells.Item(11, 3).Value = 5.8065
$ws.Cells.Item(11, 4).Value = 4.7923
$ws.Cells.Item(11, 5).Value = -23.3346
$ws.Cells.Item(12, 3).Value = 5.7749
$ws.Cells.Item(12, 4).Value = 6.0594
$ws.Cells.Item(12, 5).Value = 7.7995
$ws.Cells.Item(13, 3).Value = 5.726
$ws.Cells.Item(13, 4).Value = 11.6166
$ws.Cells.Item(13, 5).Value = 13.7152
$ws.Cells.Item(14, 3).Value = 5.6949
$ws.Cells.Item(14, 4).Value = 9.0436
$ws.Cells.Item(14, 5).Value = 15.2129
$ws.Cells.Item(17, 3).Value = 5.112
$ws.Cells.Item(17, 4).Value = 7.4628
$ws.Cells.Item(17, 5).Value = 6.4336
$ws.Cells.Item(18, 3).Value = 5.0794
$ws.Cells.Item(18, 4).Value = 3.6907
$ws.Cells.Item(18, 5).Value = 9.0009
$ws.Cells.Item(22, 3).Value = 4.9892
$ws.Cells.Item(22, 4).Value = 12.5016
$ws.Cells.Item(22, 5).Value = 10.774
$ws.Cells.Item(23, 3).Value = 4.8976
$ws.Cells.Item(23, 4).Value = 10.803
$ws.Cells.Item(23, 5).Value = 27.853
$ws.Cells.Item(25, 3).Value = 4.6905
$ws.Cells.Item(25, 4).Value = 11.241
$ws.Cells.Item(25, 5).Value = 7.7375
$ws.Cells.Item(26, 3).Value = 4.4872
$ws.Cells.Item(26, 4).Value = 17.7281
$ws.Cells.Item(26, 5).Value = 17.8159
$ws.Cells.Item(27, 3).Value = 4.4177
$ws.Cells.Item(27, 4).Value = 4.3129
$ws.Cells.Item(27, 5).Value = 5.9063
$ws.Cells.Item(29, 3).Value = 4.3728
$ws.Cells.Item(29, 4).Value = 1.6303
$ws.Cells.Item(29, 5).Value = 3.1033
$ws.Cells.Item(31, 2).Value = "SKYGOLD"
$ws.Cells.Item(31, 3).Value = 4.234
$ws.Cells.Item(31, 4).Value = -0.3419
$ws.Cells.Item(31, 5).Value = 38.4046
$ws.Cells.Item(32, 2).Value = "REDTAPE"
$ws.Cells.Item(32, 3).Value = 4.2289
$ws.Cells.Item(32, 4).Value = 4.1278
$ws.Cells.Item(32, 5).Value = -2.8416
$ws.Cells.Item(33, 2).Value = "MTARTECH"
$ws.Cells.Item(33, 3).Value = 4.1224
$ws.Cells.Item(33, 4).Value = 8.2682
$ws.Cells.Item(33, 5).Value = 32.3323
$ws.Cells.Item(34, 2).Value = "BAJAJHCARE"
$ws.Cells.Item(34, 3).Value = 4.0627
$ws.Cells.Item(34, 4).Value = 4.6097
$ws.Cells.Item(34, 5).Value = -1.6673
$ws.Cells.Item(35, 2).Value = "CENTRUM"
$ws.Cells.Item(35, 3).Value = 3.9394
$ws.Cells.Item(35, 4).Value = 2.5411
$ws.Cells.Item(35, 5).Value = 1.8711
$ws.Cells.Item(37, 3).Value = 3.8673
$ws.Cells.Item(37, 4).Value = 12.2266
$ws.Cells.Item(37, 5).Value = 11.7616
$ws.Cells.Item(38, 3).Value = 3.7113
$ws.Cells.Item(38, 4).Value = 11.0109
$ws.Cells.Item(38, 5).Value = 3.6047
$ws.Cells.Item(39, 3).Value = 3.6256
$ws.Cells.Item(39, 4).Value = 13.688
$ws.Cells.Item(39, 5).Value = 24.407
$ws.Cells.Item(40, 3).Value = 3.4968
$ws.Cells.Item(40, 4).Value = 0.4399
$ws.Cells.Item(40, 5).Value = -0.8222
$ws.Cells.Item(41, 2).Value = "ALICON"
$ws.Cells.Item(41, 3).Value = 3.468
$ws.Cells.Item(41, 4).Value = 9.601699999999999
$ws.Cells.Item(41, 5).Value = 15.0459
$ws.Cells.Item(42, 2).Value = "RSYSTEMS"
$ws.Cells.Item(42, 3).Value = 3.4387
$ws.Cells.Item(42, 4).Value = 4.5998
$ws.Cells.Item(42, 5).Value = 6.9486
$ws.Cells.Item(43, 2).Value = "PSPPROJECT"
$ws.Cells.Item(43, 3).Value = 3.3259
$ws.Cells.Item(43, 4).Value = 17.4617
$ws.Cells.Item(43, 5).Value = 23.9197
$ws.Cells.Item(44, 2).Value = "CENTUM"
$ws.Cells.Item(44, 3).Value = 3.2372
$ws.Cells.Item(44, 4).Value = 4.0991
$ws.Cells.Item(44, 5).Value = -1.4201
$ws.Cells.Item(45, 2).Value = "BGRENERGY"
$ws.Cells.Item(45, 3).Value = 3.2177
$ws.Cells.Item(45, 4).Value = -6.0894
$ws.Cells.Item(45, 5).Value = 74.8323
$ws.Cells.Item(46, 2).Value = "SHRINGARMS"
$ws.Cells.Item(46, 3).Value = 3.2134
$ws.Cells.Item(46, 4).Value = 4.4204
$ws.Cells.Item(46, 5).Value = 24.5912
$ws.Cells.Item(47, 2).Value = "IVALUE"
$ws.Cells.Item(47, 3).Value = 3.0936
$ws.Cells.Item(47, 4).Value = 6.3437
$ws.Cells.Item(47, 5).Value = -0.8602
$ws.Cells.Item(48, 2).Value = "OIL"
$ws.Cells.Item(48, 3).Value = 3.044
$ws.Cells.Item(48, 4).Value = 3.2896
$ws.Cells.Item(48, 5).Value = 4.7124
$ws.Cells.Item(49, 2).Value = "ASHOKA"
$ws.Cells.Item(49, 3).Value = 3.023
$ws.Cells.Item(49, 4).Value = 4.5462
$ws.Cells.Item(49, 5).Value = 7.1777
$ws.Cells.Item(50, 2).Value = "GMMPFAUDLR"
$ws.Cells.Item(50, 3).Value = 3.018
$ws.Cells.Item(50, 4).Value = 7.4826
$ws.Cells.Item(50, 5).Value = 19.8264
$ws.Cells.Item(51, 2).Value = "VSTIND"
$ws.Cells.Item(51, 3).Value = 2.9975
$ws.Cells.Item(51, 4).Value = 3.4979
$ws.Cells.Item(51, 5).Value = 3.0373
$ws.Cells.Item(52, 2).Value = "GANESHCP"
$ws.Cells.Item(52, 3).Value = 2.9745
$ws.Cells.Item(52, 4).Value = 2.4485
$ws.Cells.Item(52, 5).Value = 1.9993
$ws.Cells.Item(53, 2).Value = "NEULANDLAB"
$ws.Cells.Item(53, 3).Value = 2.9664
$ws.Cells.Item(53, 4).Value = -1.3986
$ws.Cells.Item(53, 5).Value = 8.6135
$ws.Cells.Item(54, 2).Value = "BLISSGVS"
$ws.Cells.Item(54, 3).Value = 2.8938
$ws.Cells.Item(54, 4).Value = 2.2482
$ws.Cells.Item(54, 5).Value = 2.6171
$ws.Cells.Item(55, 2).Value = "SUNDROP"
$ws.Cells.Item(55, 3).Value = 2.8612
$ws.Cells.Item(55, 4).Value = 2.7541
$ws.Cells.Item(55, 5).Value = 0.8563
$ws.Cells.Item(56, 2).Value = "SPANDANA"
$ws.Cells.Item(56, 3).Value = 2.8304
$ws.Cells.Item(56, 4).Value = 4.3276
$ws.Cells.Item(56, 5).Value = 2.9412
$ws.Cells.Item(57, 2).Value = "MFSL"
$ws.Cells.Item(57, 3).Value = 2.7851
$ws.Cells.Item(57, 4).Value = 2.8393
$ws.Cells.Item(57, 5).Value = -0.9391
$ws.Cells.Item(58, 2).Value = "BPCL"
$ws.Cells.Item(58, 3).Value = 2.7722
$ws.Cells.Item(58, 4).Value = 8.2615
$ws.Cells.Item(58, 5).Value = 5.329
$ws.Cells.Item(59, 2).Value = "DBCORP"
$ws.Cells.Item(59, 3).Value = 2.7678
$ws.Cells.Item(59, 4).Value = 5.4075
$ws.Cells.Item(59, 5).Value = 1.3556
$ws.Cells.Item(60, 3).Value = 2.7504
$ws.Cells.Item(60, 4).Value = 0.1788
$ws.Cells.Item(60, 5).Value = 2.1596
$ws.Cells.Item(61, 2).Value = "CARYSIL"
$ws.Cells.Item(61, 3).Value = 2.7273
$ws.Cells.Item(61, 4).Value = 2.2113
$ws.Cells.Item(61, 5).Value = 11.1045
$ws.Cells.Item(62, 2).Value = "AHLUCONT"
$ws.Cells.Item(62, 3).Value = 2.6699
$ws.Cells.Item(62, 4).Value = 1.4817
$ws.Cells.Item(62, 5).Value = -5.7174
$ws.Cells.Item(63, 3).Value = 2.5734
$ws.Cells.Item(63, 4).Value = 2.4525
$ws.Cells.Item(63, 5).Value = 2.1906
$ws.Cells.Item(65, 2).Value = "IIFL"
$ws.Cells.Item(65, 3).Value = 2.543
$ws.Cells.Item(65, 4).Value = 9.4937
$ws.Cells.Item(65, 5).Value = 18.6767
$ws.Cells.Item(66, 2).Value = "INOXGREEN"
$ws.Cells.Item(66, 3).Value = 2.5251
$ws.Cells.Item(66, 4).Value = 10.5185
$ws.Cells.Item(66, 5).Value = 33.8521
$ws.Cells.Item(67, 2).Value = "GRMOVER"
$ws.Cells.Item(67, 3).Value = 2.5115
$ws.Cells.Item(67, 4).Value = 2.6851
$ws.Cells.Item(67, 5).Value = 18.2946
$ws.Cells.Item(68, 2).Value = "BEML"
$ws.Cells.Item(68, 3).Value = 2.4348
$ws.Cells.Item(68, 4).Value = -0.4642
$ws.Cells.Item(68, 5).Value = 5.8388
$ws.Cells.Item(69, 2).Value = "FCL"
$ws.Cells.Item(69, 3).Value = 2.4256
$ws.Cells.Item(69, 4).Value = -0.2538
$ws.Cells.Item(69, 5).Value = 2.4051
$ws.Cells.Item(70, 2).Value = "FISCHER"
$ws.Cells.Item(70, 3).Value = 2.4014
$ws.Cells.Item(70, 4).Value = 13.3916
$ws.Cells.Item(70, 5).Value = 6.4144
$ws.Cells.Item(72, 2).Value = "SDBL"
$ws.Cells.Item(72, 3).Value = 2.3826
$ws.Cells.Item(72, 4).Value = 0.9304
$ws.Cells.Item(72, 5).Value = 6.5105
$ws.Cells.Item(74, 2).Value = "POWERINDIA"
$ws.Cells.Item(74, 3).Value = 2.3167
$ws.Cells.Item(74, 4).Value = 6.6798
$ws.Cells.Item(74, 5).Value = -0.6333
$ws.Cells.Item(75, 2).Value = "JKTYRE"
$ws.Cells.Item(75, 3).Value = 2.2869
$ws.Cells.Item(75, 4).Value = 5.3162
$ws.Cells.Item(75, 5).Value = 21.3057
$ws.Cells.Item(76, 2).Value = "DBL"
$ws.Cells.Item(76, 3).Value = 2.2687
$ws.Cells.Item(76, 4).Value = 3.3501
$ws.Cells.Item(76, 5).Value = 4.4214

$ws = $wb.Worksheets.Item("Top Losers")
$ws.Cells.Item(2, 3).Value = -17.1517
$ws.Cells.Item(2, 4).Value = -15.8693
$ws.Cells.Item(2, 5).Value = 1.4517
$ws.Cells.Item(3, 3).Value = -9.858700000000001
$ws.Cells.Item(3, 4).Value = -6.4334
$ws.Cells.Item(3, 5).Value = 5.7605
$ws.Cells.Item(4, 3).Value = -6.6239
$ws.Cells.Item(4, 4).Value = -9.147600000000001
$ws.Cells.Item(4, 5).Value = 7.5031
$ws.Cells.Item(5, 3).Value = -5.1613
$ws.Cells.Item(5, 4).Value = -3.5062
$ws.Cells.Item(5, 5).Value = 19.4964
$ws.Cells.Item(11, 2).Value = "RAJRATAN"
$ws.Cells.Item(11, 3).Value = -4.475
$ws.Cells.Item(11, 4).Value = -3.6925
$ws.Cells.Item(11, 5).Value = 21.1285
$ws.Cells.Item(12, 2).Value = "LICHSGFIN"
$ws.Cells.Item(12, 3).Value = -4.3727
$ws.Cells.Item(12, 4).Value = -2.0708
$ws.Cells.Item(12, 5).Value = 0.4425
$ws.Cells.Item(13, 2).Value = "SOUTHBANK"
$ws.Cells.Item(13, 3).Value = -4.0655
$ws.Cells.Item(13, 4).Value = -1.5482
$ws.Cells.Item(13, 5).Value = 29.6028
$ws.Cells.Item(14, 2).Value = "NSLNISP"
$ws.Cells.Item(14, 3).Value = -3.9716
$ws.Cells.Item(14, 4).Value = -2.7197
$ws.Cells.Item(14, 5).Value = -3.5221
$ws.Cells.Item(16, 2).Value = "RAMCOSYS"
$ws.Cells.Item(16, 3).Value = -3.779
$ws.Cells.Item(16, 4).Value = 5.9898
$ws.Cells.Item(16, 5).Value = 24.4989
$ws.Cells.Item(17, 2).Value = "DRREDDY"
$ws.Cells.Item(17, 3).Value = -3.7093
$ws.Cells.Item(17, 4).Value = -6.1624
$ws.Cells.Item(17, 5).Value = -1.569
$ws.Cells.Item(18, 2).Value = "INDUSTOWER"
$ws.Cells.Item(18, 3).Value = -3.6347
$ws.Cells.Item(18, 4).Value = 1.5627
$ws.Cells.Item(18, 5).Value = 7.0866
$ws.Cells.Item(19, 2).Value = "YATRA"
$ws.Cells.Item(19, 3).Value = -3.6227
$ws.Cells.Item(19, 4).Value = -6.365
$ws.Cells.Item(19, 5).Value = 3.4814
$ws.Cells.Item(20, 2).Value = "SARDAEN"
$ws.Cells.Item(20, 3).Value = -3.5967
$ws.Cells.Item(20, 4).Value = -0.0835
$ws.Cells.Item(20, 5).Value = -0.1205
$ws.Cells.Item(21, 2).Value = "IDEAFORGE"
$ws.Cells.Item(21, 3).Value = -3.4845
$ws.Cells.Item(21, 4).Value = -2.6855
$ws.Cells.Item(21, 5).Value = -4.4224
$ws.Cells.Item(23, 3).Value = -3.2423
$ws.Cells.Item(23, 4).Value = -2.1829
$ws.Cells.Item(23, 5).Value = 3.3408
$ws.Cells.Item(24, 2).Value = "UBL"
$ws.Cells.Item(24, 3).Value = -3.2006
$ws.Cells.Item(24, 4).Value = -2.6078
$ws.Cells.Item(24, 5).Value = -1.0773
$ws.Cells.Item(25, 2).Value = "APOLLOPIPE"
$ws.Cells.Item(25, 3).Value = -3.1946
$ws.Cells.Item(25, 4).Value = -4.8827
$ws.Cells.Item(25, 5).Value = -9.9099
$ws.Cells.Item(26, 2).Value = "VGUARD"
$ws.Cells.Item(26, 3).Value = -3.1891
$ws.Cells.Item(26, 4).Value = -0.497
$ws.Cells.Item(26, 5).Value = -1.3715
$ws.Cells.Item(27, 2).Value = "GOKULAGRO"
$ws.Cells.Item(27, 3).Value = -3.0884
$ws.Cells.Item(27, 4).Value = 4.4698
$ws.Cells.Item(27, 5).Value = -13.9105
$ws.Cells.Item(28, 3).Value = -3.0524
$ws.Cells.Item(28, 4).Value = 12.2254
$ws.Cells.Item(28, 5).Value = 29.4513
$ws.Cells.Item(29, 3).Value = -3.0452
$ws.Cells.Item(29, 4).Value = 5.2603
$ws.Cells.Item(29, 5).Value = 1.3311
$ws.Cells.Item(30, 2).Value = "RAYMONDREL"
$ws.Cells.Item(30, 3).Value = -3.0154
$ws.Cells.Item(30, 4).Value = -4.5933
$ws.Cells.Item(30, 5).Value = 9.9215
$ws.Cells.Item(31, 3).Value = -2.9302
$ws.Cells.Item(31, 4).Value = -2.7305
$ws.Cells.Item(31, 5).Value = 5.7472
$ws.Cells.Item(32, 2).Value = "TVSELECT"
$ws.Cells.Item(32, 3).Value = -2.9173
$ws.Cells.Item(32, 4).Value = -3.8627
$ws.Cells.Item(32, 5).Value = -5.8267
$ws.Cells.Item(33, 2).Value = "DREDGECORP"
$ws.Cells.Item(33, 3).Value = -2.9108
$ws.Cells.Item(33, 4).Value = 18.2183
$ws.Cells.Item(33, 5).Value = 18.9911
$ws.Cells.Item(34, 2).Value = "EPACKPEB"
$ws.Cells.Item(34, 3).Value = -2.8688
$ws.Cells.Item(34, 4).Value = -2.8037
$ws.Cells.Item(34, 5).Value = "N/A"
$ws.Cells.Item(35, 2).Value = "BHARATWIRE"
$ws.Cells.Item(35, 3).Value = -2.8635
$ws.Cells.Item(35, 4).Value = 19.3163
$ws.Cells.Item(35, 5).Value = 20.3501
$ws.Cells.Item(36, 2).Value = "FILATEX"
$ws.Cells.Item(36, 3).Value = -2.8571
$ws.Cells.Item(36, 4).Value = 7.1233
$ws.Cells.Item(36, 5).Value = 22.4027
$ws.Cells.Item(37, 2).Value = "FABTECH"
$ws.Cells.Item(37, 3).Value = -2.7992
$ws.Cells.Item(37, 4).Value = 12.3866
$ws.Cells.Item(37, 5).Value = "N/A"
$ws.Cells.Item(38, 2).Value = "BCG"
$ws.Cells.Item(38, 3).Value = -2.7027
$ws.Cells.Item(38, 4).Value = 2.3186
$ws.Cells.Item(38, 5).Value = -1.5119
$ws.Cells.Item(39, 2).Value = "JINDALPHOT"
$ws.Cells.Item(39, 3).Value = -2.6444
$ws.Cells.Item(39, 4).Value = -2.8597
$ws.Cells.Item(39, 5).Value = 19.7311
$ws.Cells.Item(40, 3).Value = -2.6358
$ws.Cells.Item(40, 4).Value = 11.928
$ws.Cells.Item(40, 5).Value = 10.2837
$ws.Cells.Item(41, 2).Value = "BIL"
$ws.Cells.Item(41, 3).Value = -2.612
$ws.Cells.Item(41, 4).Value = 6.2719
$ws.Cells.Item(41, 5).Value = -2.9239
$ws.Cells.Item(42, 2).Value = "SANDHAR"
$ws.Cells.Item(42, 3).Value = -2.5407
$ws.Cells.Item(42, 4).Value = 1.1797
$ws.Cells.Item(42, 5).Value = 18.5681
$ws.Cells.Item(43, 2).Value = "UTIAMC"
$ws.Cells.Item(43, 3).Value = -2.5219
$ws.Cells.Item(43, 4).Value = -7.1402
$ws.Cells.Item(43, 5).Value = -4.5625
$ws.Cells.Item(44, 2).Value = "ARIHANTCAP"
$ws.Cells.Item(44, 3).Value = -2.4864
$ws.Cells.Item(44, 4).Value = 4.4628
$ws.Cells.Item(44, 5).Value = -4.3442
$ws.Cells.Item(45, 2).Value = "MANGCHEFER"
$ws.Cells.Item(45, 3).Value = -2.4673
$ws.Cells.Item(45, 4).Value = -2.9058
$ws.Cells.Item(45, 5).Value = -6.7818
$ws.Cells.Item(46, 2).Value = "SURAJEST"
$ws.Cells.Item(46, 3).Value = -2.457
$ws.Cells.Item(46, 4).Value = 6.5352
$ws.Cells.Item(46, 5).Value = 4.5316
$ws.Cells.Item(48, 2).Value = "HFCL"
$ws.Cells.Item(48, 3).Value = -2.3979
$ws.Cells.Item(48, 4).Value = -3.1387
$ws.Cells.Item(48, 5).Value = 3.3205
$ws.Cells.Item(49, 2).Value = "STYL"
$ws.Cells.Item(49, 3).Value = -2.3902
$ws.Cells.Item(49, 4).Value = -5.58
$ws.Cells.Item(49, 5).Value = -11.1678
$ws.Cells.Item(50, 2).Value = "IDBI"
$ws.Cells.Item(50, 3).Value = -2.3845
$ws.Cells.Item(50, 4).Value = 5.9087
$ws.Cells.Item(50, 5).Value = 8.8165
$ws.Cells.Item(51, 2).Value = "GRWRHITECH"
$ws.Cells.Item(51, 3).Value = -2.3845
$ws.Cells.Item(51, 4).Value = -5.9774
$ws.Cells.Item(51, 5).Value = 18.7755
$ws.Cells.Item(52, 2).Value = "CAMLINFINE"
$ws.Cells.Item(52, 3).Value = -2.3694
$ws.Cells.Item(52, 4).Value = 0.4189
$ws.Cells.Item(52, 5).Value = 0.6813
$ws.Cells.Item(53, 2).Value = "POCL"
$ws.Cells.Item(53, 3).Value = -2.352
$ws.Cells.Item(53, 4).Value = 2.8188
$ws.Cells.Item(53, 5).Value = 23.37
$ws.Cells.Item(54, 2).Value = "TBOTEK"
$ws.Cells.Item(54, 3).Value = -2.3298
$ws.Cells.Item(54, 4).Value = -5.8198
$ws.Cells.Item(54, 5).Value = -1.3179
$ws.Cells.Item(55, 3).Value = -2.3188
$ws.Cells.Item(55, 4).Value = -5.0901
$ws.Cells.Item(55, 5).Value = 1.8446
$ws.Cells.Item(56, 3).Value = -2.2829
$ws.Cells.Item(56, 4).Value = -1.1954
$ws.Cells.Item(56, 5).Value = -15.1504
$ws.Cells.Item(57, 2).Value = "NEWGEN"
$ws.Cells.Item(57, 3).Value = -2.2747
$ws.Cells.Item(57, 4).Value = 9.000299999999999
$ws.Cells.Item(57, 5).Value = 9.373900000000001
$ws.Cells.Item(58, 2).Value = "MSPL"
$ws.Cells.Item(58, 3).Value = -2.274
$ws.Cells.Item(58, 4).Value = -1.1209
$ws.Cells.Item(58, 5).Value = -8.164400000000001
$ws.Cells.Item(59, 2).Value = "VIPIND"
$ws.Cells.Item(59, 3).Value = -2.2657
$ws.Cells.Item(59, 4).Value = -3.9944
$ws.Cells.Item(59, 5).Value = -1.4785
$ws.Cells.Item(60, 2).Value = "THYROCARE"
$ws.Cells.Item(60, 4).Value = 4.2572
$ws.Cells.Item(60, 5).Value = 13.7343
$ws.Cells.Item(61, 2).Value = "VBL"
$ws.Cells.Item(61, 3).Value = -2.2606
$ws.Cells.Item(61, 4).Value = 4.9523
$ws.Cells.Item(61, 5).Value = 9.139099999999999
$ws.Cells.Item(63, 2).Value = "HCG"
$ws.Cells.Item(63, 3).Value = -2.2285
$ws.Cells.Item(63, 4).Value = -0.0132
$ws.Cells.Item(63, 5).Value = 17.9854
$ws.Cells.Item(64, 2).Value = "CGCL"
$ws.Cells.Item(64, 3).Value = -2.224
$ws.Cells.Item(64, 4).Value = -0.2586
$ws.Cells.Item(64, 5).Value = 10.3953
$ws.Cells.Item(65, 2).Value = "TTKPRESTIG"
$ws.Cells.Item(65, 3).Value = -2.2182
$ws.Cells.Item(65, 4).Value = 5.6055
$ws.Cells.Item(65, 5).Value = 7.2182
$ws.Cells.Item(66, 2).Value = "LXCHEM"
$ws.Cells.Item(66, 3).Value = -2.2147
$ws.Cells.Item(66, 4).Value = -2.6554
$ws.Cells.Item(66, 5).Value = -3.8333
$ws.Cells.Item(67, 2).Value = "SSWL"
$ws.Cells.Item(67, 3).Value = -2.1903
$ws.Cells.Item(67, 4).Value = 4.0099
$ws.Cells.Item(67, 5).Value = 1.1144
$ws.Cells.Item(68, 2).Value = "GABRIEL"
$ws.Cells.Item(68, 3).Value = -2.164
$ws.Cells.Item(68, 4).Value = 1.9725
$ws.Cells.Item(68, 5).Value = 6.7716
$ws.Cells.Item(69, 2).Value = "KFINTECH"
$ws.Cells.Item(69, 3).Value = -2.1341
$ws.Cells.Item(69, 4).Value = -3.8455
$ws.Cells.Item(69, 5).Value = 5.087
$ws.Cells.Item(70, 2).Value = "NUVAMA"
$ws.Cells.Item(70, 3).Value = -2.1271
$ws.Cells.Item(70, 4).Value = 0.3484
$ws.Cells.Item(70, 5).Value = 14.2109
$ws.Cells.Item(71, 2).Value = "SURYAROSNI"
$ws.Cells.Item(71, 3).Value = -2.119
$ws.Cells.Item(71, 4).Value = 9.025700000000001
$ws.Cells.Item(71, 5).Value = 0.8383
$ws.Cells.Item(72, 2).Value = "TMB"
$ws.Cells.Item(72, 3).Value = -2.1042
$ws.Cells.Item(72, 4).Value = 7.4864
$ws.Cells.Item(72, 5).Value = 14.67
$ws.Cells.Item(73, 2).Value = "GARUDA"
$ws.Cells.Item(73, 3).Value = -2.1
$ws.Cells.Item(73, 4).Value = -8.865600000000001
$ws.Cells.Item(73, 5).Value = 7.4675
$ws.Cells.Item(74, 2).Value = "SOLARWORLD"
$ws.Cells.Item(74, 3).Value = -2.0981
$ws.Cells.Item(74, 4).Value = 6.7332
$ws.Cells.Item(74, 5).Value = 2.4111
$ws.Cells.Item(75, 2).Value = "CGPOWER"
$ws.Cells.Item(75, 3).Value = -2.0972
$ws.Cells.Item(75, 4).Value = 1.2503
$ws.Cells.Item(75, 5).Value = -1.0864
$ws.Cells.Item(76, 2).Value = "WAAREERTL"
$ws.Cells.Item(76, 3).Value = -2.0845
$ws.Cells.Item(76, 4).Value = 1.3004
$ws.Cells.Item(76, 5).Value = 19.6756

$ws = $wb.Worksheets.Item("1 Month Performance")
$ws.Cells.Item(5, 3).Value = 65.33710000000001
$ws.Cells.Item(6, 3).Value = 62.1345
$ws.Cells.Item(7, 3).Value = 54.7315
$ws.Cells.Item(8, 3).Value = 53.6034
$ws.Cells.Item(9, 3).Value = 52.2629
$ws.Cells.Item(10, 3).Value = 45.9698
$ws.Cells.Item(11, 3).Value = 41.7496
$ws.Cells.Item(14, 3).Value = 39.1081
$ws.Cells.Item(16, 3).Value = 37.2059
$ws.Cells.Item(17, 3).Value = 36.5964
$ws.Cells.Item(19, 3).Value = 34.8702
$ws.Cells.Item(20, 3).Value = 34.4425
$ws.Cells.Item(21, 3).Value = 33.2997
$ws.Cells.Item(22, 3).Value = 33.0664
$ws.Cells.Item(24, 3).Value = 30.0014
$ws.Cells.Item(25, 3).Value = 29.9169
$ws.Cells.Item(26, 2).Value = "TARACHAND"
$ws.Cells.Item(26, 3).Value = 29.3749
$ws.Cells.Item(27, 2).Value = "MRPL"
$ws.Cells.Item(27, 3).Value = 29.1937
$ws.Cells.Item(28, 3).Value = 29.0193
$ws.Cells.Item(30, 2).Value = "INDORAMA"
$ws.Cells.Item(30, 3).Value = 28.3343
$ws.Cells.Item(33, 2).Value = "SKYGOLD"
$ws.Cells.Item(33, 3).Value = 25.7463
$ws.Cells.Item(34, 2).Value = "HATSUN"
$ws.Cells.Item(34, 3).Value = 25.6297
$ws.Cells.Item(35, 2).Value = "MARINE"
$ws.Cells.Item(35, 3).Value = 25.5216
$ws.Cells.Item(37, 2).Value = "TDPOWERSYS"
$ws.Cells.Item(37, 3).Value = 25.3608
$ws.Cells.Item(38, 2).Value = "ATHERENERG"
$ws.Cells.Item(38, 3).Value = 24.4472
$ws.Cells.Item(39, 2).Value = "SAGILITY"
$ws.Cells.Item(39, 3).Value = 24.4132
$ws.Cells.Item(40, 2).Value = "UNIPARTS"
$ws.Cells.Item(40, 3).Value = 24.3817
$ws.Cells.Item(42, 3).Value = 24.2949
$ws.Cells.Item(43, 3).Value = 23.9668
$ws.Cells.Item(44, 2).Value = "AUBANK"
$ws.Cells.Item(44, 3).Value = 23.615
$ws.Cells.Item(45, 2).Value = "RAMCOSYS"
$ws.Cells.Item(45, 3).Value = 23.4863
$ws.Cells.Item(46, 3).Value = 23.2208
$ws.Cells.Item(48, 3).Value = 22.5792
$ws.Cells.Item(49, 3).Value = 22.4806
$ws.Cells.Item(51, 3).Value = 22.0539
$ws.Cells.Item(52, 3).Value = 21.6253
$ws.Cells.Item(53, 3).Value = 21.3115
$ws.Cells.Item(54, 2).Value = "GRMOVER"
$ws.Cells.Item(54, 3).Value = 21.2857
$ws.Cells.Item(55, 2).Value = "SURYODAY"
$ws.Cells.Item(55, 3).Value = 21.2719
$ws.Cells.Item(56, 2).Value = "CPEDU"
$ws.Cells.Item(56, 3).Value = 21.0396
$ws.Cells.Item(59, 2).Value = "BHAGERIA"
$ws.Cells.Item(59, 3).Value = 20.6198
$ws.Cells.Item(61, 2).Value = "INDRAMEDCO"
$ws.Cells.Item(61, 3).Value = 20.259
$ws.Cells.Item(62, 3).Value = 20.0725
$ws.Cells.Item(63, 2).Value = "IIFL"
$ws.Cells.Item(63, 3).Value = 20.0582
$ws.Cells.Item(64, 2).Value = "BHARATWIRE"
$ws.Cells.Item(64, 3).Value = 19.8379
$ws.Cells.Item(65, 2).Value = "FEDERALBNK"
$ws.Cells.Item(65, 3).Value = 19.6324
$ws.Cells.Item(66, 2).Value = "SHRIRAMFIN"
$ws.Cells.Item(66, 3).Value = 19.6247
$ws.Cells.Item(67, 2).Value = "BANKINDIA"
$ws.Cells.Item(67, 3).Value = 19.4589
$ws.Cells.Item(70, 3).Value = 19.04
$ws.Cells.Item(71, 3).Value = 18.8832
$ws.Cells.Item(74, 2).Value = "THOMASCOTT"
$ws.Cells.Item(74, 3).Value = 18.3194
$ws.Cells.Item(75, 2).Value = "REPRO"
$ws.Cells.Item(75, 3).Value = 18.1879

$ws = $wb.Worksheets.Item("distance from Dma50")
$ws.Cells.Item(2, 3).Value = 10.0381
$ws.Cells.Item(3, 3).Value = 7.4941
$ws.Cells.Item(4, 3).Value = 6.3364
$ws.Cells.Item(5, 3).Value = 5.3117
$ws.Cells.Item(6, 3).Value = 5.2626
$ws.Cells.Item(7, 3).Value = 5.0566
$ws.Cells.Item(8, 3).Value = 4.4785
$ws.Cells.Item(9, 3).Value = 4.3967
$ws.Cells.Item(10, 3).Value = 3.8748
$ws.Cells.Item(11, 3).Value = 3.5988
$ws.Cells.Item(12, 3).Value = 3.3888
$ws.Cells.Item(13, 3).Value = 3.3626
$ws.Cells.Item(14, 3).Value = 3.108
$ws.Cells.Item(15, 3).Value = 3.0695
$ws.Cells.Item(16, 3).Value = 2.9927
$ws.Cells.Item(17, 3).Value = 2.8431
$ws.Cells.Item(18, 3).Value = 2.6361
$ws.Cells.Item(19, 3).Value = 2.6324
$ws.Cells.Item(20, 3).Value = 2.3918
$ws.Cells.Item(21, 3).Value = 2.3674
$ws.Cells.Item(22, 3).Value = 1.4684
$ws.Cells.Item(23, 3).Value = 1.4094
$ws.Cells.Item(24, 3).Value = 1.4036
$ws.Cells.Item(25, 3).Value = 1.2265
$ws.Cells.Item(26, 3).Value = 1.0751
$ws.Cells.Item(27, 3).Value = 0.9844000000000001
$ws.Cells.Item(28, 3).Value = 0.6417
$ws.Cells.Item(29, 3).Value = 0.28
$ws.Cells.Item(30, 3).Value = -1.9974
